$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates ---
# Most of these strings (single "." separator, e.g. "336.31") would be
# silently re-interpreted by Excel as a floating point number on plain
# assignment, losing exact digits/trailing zeros (e.g. "6.850" -> 6.85).
# To keep them as literal text we temporarily force the cell to Text
# format, assign the value, then clear the temporary format again so the
# cell keeps its original (default) style, same as in the source file.
# Values containing two "." separators (e.g. "30.552.37") are never
# parsed as numbers by Excel, so a plain assignment is sufficient there.

$ws.Range("D2").Value = "30.552.37"
$ws.Range("D3").Value = "2.109.35"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.31"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5252"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4553"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.74"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09003"
$ws.Range("D10").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.55"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "2.118.10"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.850"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.112"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001177"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.25"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.33"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.262"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "30.616.85"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "2.356.75"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.62"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.517"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.59"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.220"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1069"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.343"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.626"
$ws.Range("D34").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.52"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.888"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02608"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06819"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2314"
$ws.Range("D40").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6851"
$ws.Range("D42").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.306"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.676"
$ws.Range("D47").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000350"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.210"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.16"
$ws.Range("D51").ClearFormats()

# --- Column E ("Volume(1h)") updates ---
# These are always padded percentage strings (e.g. "  +0.36%  "); Excel
# never interprets them as numbers, so plain assignment is sufficient.

$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +1.02%  "
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +4.77%  "
$ws.Range("E9").Value = "  +4.78%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("E15").Value = "  +5.65%  "
$ws.Range("E16").Value = "  +5.20%  "
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("E36").Value = "  +4.34%  "
$ws.Range("E37").Value = "  +8.02%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +19.18%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  +1.56%  "
